$wb = $excel.ActiveWorkbook

# --- Summary sheet: updated capital / P&L / trade-count metrics ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.23   # Current Capital
$summary.Range("B4").Value = 0.24      # Total P&L $
$summary.Range("B5").Value = 0.04      # Total P&L %
$summary.Range("B6").Value = 113       # Total Trades
$summary.Range("B7").Value = 50        # Winning Trades
$summary.Range("B9").Value = 44.25     # Win Rate %

# --- Strategy Status sheet: MarketMaking row (row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.23     # Capital
$status.Range("D4").Value = 113        # Trades
$status.Range("E4").Value = 0.24       # P&L $
$status.Range("F4").Value = 0.23       # P&L %
$status.Range("G4").Value = 44.25      # Win Rate %

# --- New trade row (#113) appended to both "All Trades" and "MarketMaking" sheets ---
# Helper: writes the 17-column trade record into row 114 of a given worksheet.
# Date/time-looking text values are entered with a leading apostrophe so the
# COM layer keeps them as literal text instead of auto-converting them to
# Excel date/time serials (matching how the rest of the sheet is stored),
# then the quote-prefix formatting is cleared via Style="Normal" so no extra
# style survives on the cell.
function Set-TradeRow114($sheet) {
    $sheet.Cells.Item(114, 1).Value = 113
    $sheet.Cells.Item(114, 2).Value = "'2026-02-17"
    $sheet.Cells.Item(114, 2).Style = "Normal"
    $sheet.Cells.Item(114, 3).Value = "'09:19:45"
    $sheet.Cells.Item(114, 3).Style = "Normal"
    $sheet.Cells.Item(114, 4).Value = "MarketMaking"
    $sheet.Cells.Item(114, 5).Value = "DOWN"
    $sheet.Cells.Item(114, 6).Value = 0.74
    $sheet.Cells.Item(114, 7).Value = 0.809094
    $sheet.Cells.Item(114, 8).Value = "CLOSED"
    $sheet.Cells.Item(114, 9).Value = 9.3371
    $sheet.Cells.Item(114, 10).Value = 0.07000000000000001
    $sheet.Cells.Item(114, 11).Value = 100.23
    $sheet.Cells.Item(114, 12).Value = 0
    $sheet.Cells.Item(114, 13).Value = 0
    $sheet.Cells.Item(114, 14).Value = 0.6
    $sheet.Cells.Item(114, 15).Value = "Normal spread capture: 19600 bps"
    $sheet.Cells.Item(114, 16).Value = "early_exit"
    $sheet.Cells.Item(114, 17).Value = 0.11
}

Set-TradeRow114 $wb.Worksheets.Item("All Trades")
Set-TradeRow114 $wb.Worksheets.Item("MarketMaking")
